# Update "想去人数" (want-to-go count) figures in column F across sheets,
# mirroring the refreshed scrape output committed at 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F5").Value = 382
$wsExpo.Range("F6").Value = 366
$wsExpo.Range("F8").Value = 243
$wsExpo.Range("F9").Value = 786
$wsExpo.Range("F10").Value = 2359
$wsExpo.Range("F11").Value = 359
$wsExpo.Range("F13").Value = 236
$wsExpo.Range("F15").Value = 209
$wsExpo.Range("F16").Value = 199
$wsExpo.Range("F17").Value = 2803
$wsExpo.Range("F22").Value = 239
$wsExpo.Range("F24").Value = 278

# Sheet "本地生活" (local life)
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F4").Value = 2126

# Sheet "全部类型" (all types - aggregated view)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 2126
$wsAll.Range("F12").Value = 382
$wsAll.Range("F17").Value = 366
$wsAll.Range("F20").Value = 243
$wsAll.Range("F24").Value = 786
$wsAll.Range("F25").Value = 2359
$wsAll.Range("F26").Value = 359
$wsAll.Range("F29").Value = 236
$wsAll.Range("F31").Value = 209
$wsAll.Range("F32").Value = 199
$wsAll.Range("F41").Value = 239
$wsAll.Range("F50").Value = 278
